# Weekly fruit/vegetable price update: insert 3 new records for Espárragos
# (Mercado Mayorista Lo Valledor de Santiago) as the new rows 12-14, pushing
# the existing rows 12-18 down to rows 15-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows right before the current row 12 (inherits the
# date-format style from the row above for column D, same as the existing
# data rows).
$ws.Range("A12:A14").EntireRow.Insert()

# New row 12
$ws.Cells.Item(12, 1).Value = 6
$ws.Cells.Item(12, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(12, 3).Value = "Metropolitana"
$ws.Cells.Item(12, 4).Value = 44474
$ws.Cells.Item(12, 5).Value = 13
$ws.Cells.Item(12, 6).Value = 300000000
$ws.Cells.Item(12, 7).Value = "Espárragos"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Banquete"
$ws.Cells.Item(12, 10).Value = 780
$ws.Cells.Item(12, 11).Value = 1500
$ws.Cells.Item(12, 12).Value = 1600
$ws.Cells.Item(12, 13).Value = 1558
$ws.Cells.Item(12, 14).Value = "`$/kilo"
$ws.Cells.Item(12, 15).Value = "Provincia de Linares"
$ws.Cells.Item(12, 16).Value = 1558
$ws.Cells.Item(12, 17).Value = 1
$ws.Cells.Item(12, 18).Value = "Hortaliza"

# New row 13
$ws.Cells.Item(13, 1).Value = 6
$ws.Cells.Item(13, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(13, 3).Value = "Metropolitana"
$ws.Cells.Item(13, 4).Value = 44474
$ws.Cells.Item(13, 5).Value = 13
$ws.Cells.Item(13, 6).Value = 300000000
$ws.Cells.Item(13, 7).Value = "Espárragos"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 520
$ws.Cells.Item(13, 11).Value = 1300
$ws.Cells.Item(13, 12).Value = 1400
$ws.Cells.Item(13, 13).Value = 1348
$ws.Cells.Item(13, 14).Value = "`$/kilo"
$ws.Cells.Item(13, 15).Value = "Provincia de Linares"
$ws.Cells.Item(13, 16).Value = 1348
$ws.Cells.Item(13, 17).Value = 1
$ws.Cells.Item(13, 18).Value = "Hortaliza"

# New row 14
$ws.Cells.Item(14, 1).Value = 6
$ws.Cells.Item(14, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(14, 3).Value = "Metropolitana"
$ws.Cells.Item(14, 4).Value = 44474
$ws.Cells.Item(14, 5).Value = 13
$ws.Cells.Item(14, 6).Value = 300000000
$ws.Cells.Item(14, 7).Value = "Espárragos"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Segunda"
$ws.Cells.Item(14, 10).Value = 400
$ws.Cells.Item(14, 11).Value = 1000
$ws.Cells.Item(14, 12).Value = 1200
$ws.Cells.Item(14, 13).Value = 1100
$ws.Cells.Item(14, 14).Value = "`$/kilo"
$ws.Cells.Item(14, 15).Value = "Provincia de Linares"
$ws.Cells.Item(14, 16).Value = 1100
$ws.Cells.Item(14, 17).Value = 1
$ws.Cells.Item(14, 18).Value = "Hortaliza"

# Refresh the sheet's declared dimension to match the new used range.
$wb.Worksheets.Item(1).UsedRange | Out-Null
